# Change qty bug Fixed
#
# 1. Stocks sheet: fix Quantity of "Stock Name two" (row 3, column D) from 435 to 425.
# 2. Bills sheet: correct the balance amount of receipt XX0812128 (row 128, column E)
#    from "0.0" to "0".
# 3. Bills sheet: add the two new bills (and their clearance entry) recorded on
#    10-Dec-2020 for Aakash / Stock Name two(5).

$wb = $excel.ActiveWorkbook
$stocks = $wb.Worksheets.Item("Stocks")
$bills  = $wb.Worksheets.Item("Bills")

# Helper: force a cell to be stored as text (shared string) regardless of whether
# the value looks numeric/date-like, matching how this workbook stores every
# Bills column as text, then reset the style back to the default "Normal" style
# so no extra/quote-prefixed cell style gets introduced.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- 1. Fix the Quantity bug for "Stock Name two" ---
$stocks.Cells.Item(3, 4).Value = 425

# --- 2. Fix the Bal_Amt text for receipt XX0812128 (row 128) ---
Set-TextValue $bills.Cells.Item(128, 5) "0"

# --- 3. Append the new Bills rows 146-148 ---
Set-TextValue $bills.Cells.Item(146, 1) "10-Dec-2020 13:38"
Set-TextValue $bills.Cells.Item(146, 2) "Aakash"
Set-TextValue $bills.Cells.Item(146, 3) "364"
Set-TextValue $bills.Cells.Item(146, 4) "1000"
Set-TextValue $bills.Cells.Item(146, 5) "0"
Set-TextValue $bills.Cells.Item(146, 6) "XX1012146"
Set-TextValue $bills.Cells.Item(146, 7) "Stock Name two(5)"

Set-TextValue $bills.Cells.Item(147, 1) "10-Dec-2020 13:44"
Set-TextValue $bills.Cells.Item(147, 2) "Aakash"
Set-TextValue $bills.Cells.Item(147, 3) "364"
Set-TextValue $bills.Cells.Item(147, 4) "800.0"
Set-TextValue $bills.Cells.Item(147, 5) "0"
Set-TextValue $bills.Cells.Item(147, 6) "XX1012147"
Set-TextValue $bills.Cells.Item(147, 7) "Stock Name two(5)"

Set-TextValue $bills.Cells.Item(148, 1) "10-Dec-2020 13:44"
Set-TextValue $bills.Cells.Item(148, 2) "Aakash"
Set-TextValue $bills.Cells.Item(148, 3) "364"
Set-TextValue $bills.Cells.Item(148, 4) "1255"
Set-TextValue $bills.Cells.Item(148, 5) "0.0"
Set-TextValue $bills.Cells.Item(148, 6) "XX1012147"
Set-TextValue $bills.Cells.Item(148, 7) "Stock Name two(5),Bill Clearance 10Dec2020(1)"
